$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = "W333"
$ws.Range("I2").Value = "W333 200017758"
$ws.Range("H3").Value = "W333"
$ws.Range("I3").Value = "W333 200017758"

$ws.Range("D2").Select()
